$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Insert a new "Meta description" paragraph right after the title
#    (Heading1) paragraph at the top of the document. We use InsertXML
#    with a raw OOXML fragment (rather than InsertParagraphAfter + typed
#    text) so the paragraph gets the same "leading empty run" shape that
#    the rest of this document's body paragraphs use. InsertXML splits
#    in a new paragraph after ours, so we immediately delete that extra
#    empty paragraph.
# ---------------------------------------------------------------------
$titlePara = $d.Paragraphs.First
$insertPoint = $d.Range($titlePara.Range.End, $titlePara.Range.End)
$insertPoint.InsertXML("<pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'><pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'><pkg:xmlData><w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>: Explore the exciting world of Dragon's Fire Megaways with our detailed review. Play now for free and discover exclusive features.</w:t></w:r></w:p><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>")

$strayPara = $d.Paragraphs(3)
$strayPara.Range.Delete()

# ---------------------------------------------------------------------
# 2. Remove the duplicated bold "Play Dragon's Fire Megaways Free:
#    Review & Features" paragraph near the end of the document, and
#    replace the italic paragraph's text that follows it with the new
#    image-generation prompt text.
# ---------------------------------------------------------------------
$count = $d.Paragraphs.Count
for ($i = $count; $i -ge 1; $i--) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -eq "Play Dragon's Fire Megaways Free: Review & Features`r") {
        $p.Range.Delete()
        break
    }
}

$count = $d.Paragraphs.Count
for ($i = $count; $i -ge 1; $i--) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -eq "Explore the exciting world of Dragon's Fire Megaways with our detailed review. Play now for free and discover exclusive features.`r") {
        $pr = $p.Range
        $textOnly = $d.Range($pr.Start, $pr.End - 1)
        $textOnly.Text = "Prompt: Create an eye-catching feature image fitting the game `"Dragon's Fire Megaways`" in cartoon style. The image should feature a happy Maya warrior with glasses. For the feature image for Dragon's Fire Megaways, I suggest depicting a happy Maya warrior holding a flaming dragon egg. The warrior could be wearing glasses to bring in a modern twist to the otherwise ancient theme of the slot game. The background of the image can be a fiery orange to represent the dragon's breath, and the title of the game can be displayed in bold letters in a fiery font. The overall design of the image should be eye-catching and vibrant to capture the essence of the game."
        break
    }
}
